$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Преодолеет страх" -> two paragraphs: "Преодолеет" / " страх"
# -----------------------------------------------------------------
$target1 = "Преодолеет страх"
$rng1 = $d.Content
[void]$rng1.Find.Execute($target1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
    # Split right before the space that precedes "страх" (after the
    # 10 characters of "Преодолеет") so the break lands exactly where
    # the word boundary was, preserving the leading space on part two.
    $splitPos1 = $rng1.Start + 10
    [void]$d.Range($splitPos1, $splitPos1).InsertParagraphAfter()
}

# -----------------------------------------------------------------
# 2) "Оставшись в наших сердцах" -> three paragraphs:
#    "Оставшись" / " в наших" / "  сердцах"
#    (note: the final fragment gains an extra leading space versus
#    the single space that separated "наших" and "сердцах" before)
# -----------------------------------------------------------------
$target2 = "Оставшись в наших сердцах"
$rng2 = $d.Content
[void]$rng2.Find.Execute($target2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $base2 = $rng2.Start
    # Insert the extra space in front of "сердцах" (offset 17) first,
    # so later offsets into the still-single paragraph stay valid.
    [void]$d.Range($base2 + 17, $base2 + 17).InsertBefore(" ")
    # Split before the (now doubled) space preceding "сердцах".
    [void]$d.Range($base2 + 17, $base2 + 17).InsertParagraphAfter()
    # Split before the space preceding "в наших".
    [void]$d.Range($base2 + 9, $base2 + 9).InsertParagraphAfter()
}

Write-Output "edit applied"
